$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6825852394104004
$ws.Range("B1").Value = 0.5306469202041626
$ws.Range("C1").Value = 3.549840927124023
$ws.Range("D1").Value = 3.506358861923218
$ws.Range("E1").Value = 0.9730880856513977
